# Populates sentence-block rows 56-147 on Sheet1 (word/index/count table)
# that were appended to the "raw sentence.xlsx" dataset export.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- cell values (numbers and shared-string text) ---
$cellValues = @(
    @{Cell="D56"; Value=12}
    @{Cell="D57"; Value=13}
    @{Cell="D58"; Value=14}
    @{Cell="E58"; Value='[b''independence'']'}
    @{Cell="F58"; Value=1}
    @{Cell="D59"; Value=15}
    @{Cell="E59"; Value='[b''is'']'}
    @{Cell="F59"; Value=1}
    @{Cell="D60"; Value=16}
    @{Cell="E60"; Value='[b''essential'']'}
    @{Cell="F60"; Value=1}
    @{Cell="D61"; Value=17}
    @{Cell="E61"; Value='[b''to'']'}
    @{Cell="F61"; Value=1}
    @{Cell="D62"; Value=18}
    @{Cell="E62"; Value='[b''our'']'}
    @{Cell="F62"; Value=1}
    @{Cell="D63"; Value=19}
    @{Cell="E63"; Value='[b''legitimacy'']'}
    @{Cell="F63"; Value=1}
    @{Cell="C65"; Value=3}
    @{Cell="D65"; Value=0}
    @{Cell="D66"; Value=1}
    @{Cell="D67"; Value=2}
    @{Cell="D68"; Value=3}
    @{Cell="D69"; Value=4}
    @{Cell="E69"; Value='[b''work'']'}
    @{Cell="F69"; Value=1}
    @{Cell="D70"; Value=5}
    @{Cell="E70"; Value='[b''we'']'}
    @{Cell="F70"; Value=1}
    @{Cell="D71"; Value=6}
    @{Cell="E71"; Value='[b''help'']'}
    @{Cell="F71"; Value=1}
    @{Cell="D72"; Value=7}
    @{Cell="E72"; Value='[b''the'']'}
    @{Cell="F72"; Value=1}
    @{Cell="D73"; Value=8}
    @{Cell="E73"; Value='[b''states'']'}
    @{Cell="F73"; Value=1}
    @{Cell="D74"; Value=9}
    @{Cell="E74"; Value='[b''to'']'}
    @{Cell="F74"; Value=1}
    @{Cell="D75"; Value=10}
    @{Cell="E75"; Value='[b''implement'']'}
    @{Cell="F75"; Value=1}
    @{Cell="D76"; Value=11}
    @{Cell="E76"; Value='[b''their'']'}
    @{Cell="F76"; Value=1}
    @{Cell="D77"; Value=12}
    @{Cell="E77"; Value='[b''obligations'']'}
    @{Cell="F77"; Value=1}
    @{Cell="D78"; Value=13}
    @{Cell="E78"; Value='[b''towards'']'}
    @{Cell="F78"; Value=1}
    @{Cell="D79"; Value=14}
    @{Cell="E79"; Value='[b''all'']'}
    @{Cell="F79"; Value=1}
    @{Cell="D80"; Value=15}
    @{Cell="E80"; Value='[b''human'']'}
    @{Cell="F80"; Value=1}
    @{Cell="D81"; Value=16}
    @{Cell="E81"; Value='[b''beings'']'}
    @{Cell="F81"; Value=1}
    @{Cell="D82"; Value=17}
    @{Cell="E82"; Value='[b''and'']'}
    @{Cell="F82"; Value=1}
    @{Cell="D83"; Value=18}
    @{Cell="E83"; Value='[b''their'']'}
    @{Cell="F83"; Value=1}
    @{Cell="D84"; Value=19}
    @{Cell="E84"; Value='[b''countries'']'}
    @{Cell="F84"; Value=1}
    @{Cell="C86"; Value=4}
    @{Cell="D86"; Value=0}
    @{Cell="E86"; Value='[b''us'']'}
    @{Cell="F86"; Value=1}
    @{Cell="D87"; Value=1}
    @{Cell="E87"; Value='[b''to'']'}
    @{Cell="F87"; Value=1}
    @{Cell="D88"; Value=2}
    @{Cell="E88"; Value='[b''do'']'}
    @{Cell="F88"; Value=1}
    @{Cell="D89"; Value=3}
    @{Cell="E89"; Value='[b''that'']'}
    @{Cell="F89"; Value=1}
    @{Cell="D90"; Value=4}
    @{Cell="E90"; Value='[b''effectively'']'}
    @{Cell="F90"; Value=1}
    @{Cell="D91"; Value=5}
    @{Cell="E91"; Value='[b''we'']'}
    @{Cell="F91"; Value=1}
    @{Cell="D92"; Value=6}
    @{Cell="E92"; Value='[b''need'']'}
    @{Cell="F92"; Value=1}
    @{Cell="D93"; Value=7}
    @{Cell="E93"; Value='[b''to'']'}
    @{Cell="F93"; Value=1}
    @{Cell="D94"; Value=8}
    @{Cell="E94"; Value='[b''be'']'}
    @{Cell="F94"; Value=1}
    @{Cell="D95"; Value=9}
    @{Cell="E95"; Value='[b''independent'']'}
    @{Cell="F95"; Value=1}
    @{Cell="D96"; Value=10}
    @{Cell="E96"; Value='[b''and'']'}
    @{Cell="F96"; Value=1}
    @{Cell="D97"; Value=11}
    @{Cell="E97"; Value='[b''we'']'}
    @{Cell="F97"; Value=1}
    @{Cell="D98"; Value=12}
    @{Cell="E98"; Value='[b''need'']'}
    @{Cell="F98"; Value=1}
    @{Cell="D99"; Value=13}
    @{Cell="E99"; Value='[b''to'']'}
    @{Cell="F99"; Value=1}
    @{Cell="D100"; Value=14}
    @{Cell="E100"; Value='[b''be'']'}
    @{Cell="F100"; Value=1}
    @{Cell="D101"; Value=15}
    @{Cell="E101"; Value='[b''perceived'']'}
    @{Cell="F101"; Value=1}
    @{Cell="D102"; Value=16}
    @{Cell="E102"; Value='[b''to'']'}
    @{Cell="F102"; Value=1}
    @{Cell="D103"; Value=17}
    @{Cell="E103"; Value='[b''be'']'}
    @{Cell="F103"; Value=1}
    @{Cell="D104"; Value=18}
    @{Cell="E104"; Value='[b''independent'']'}
    @{Cell="F104"; Value=1}
    @{Cell="D105"; Value=19}
    @{Cell="E105"; Value='[b''work'']'}
    @{Cell="F105"; Value=1}
    @{Cell="C107"; Value=5}
    @{Cell="D107"; Value=0}
    @{Cell="E107"; Value='[b''rights'']'}
    @{Cell="F107"; Value=1}
    @{Cell="D108"; Value=1}
    @{Cell="E108"; Value='[b''institutions'']'}
    @{Cell="F108"; Value=1}
    @{Cell="D109"; Value=2}
    @{Cell="E109"; Value='[b''was'']'}
    @{Cell="F109"; Value=1}
    @{Cell="D110"; Value=3}
    @{Cell="E110"; Value='[b''conceived'']'}
    @{Cell="F110"; Value=1}
    @{Cell="D111"; Value=4}
    @{Cell="E111"; Value='[b''as'']'}
    @{Cell="F111"; Value=1}
    @{Cell="D112"; Value=5}
    @{Cell="E112"; Value='[b''part'']'}
    @{Cell="F112"; Value=1}
    @{Cell="D113"; Value=6}
    @{Cell="E113"; Value='[b''of'']'}
    @{Cell="F113"; Value=1}
    @{Cell="D114"; Value=7}
    @{Cell="E114"; Value='[b''the'']'}
    @{Cell="F114"; Value=1}
    @{Cell="D115"; Value=8}
    @{Cell="E115"; Value='[b''copenhagen'']'}
    @{Cell="F115"; Value=1}
    @{Cell="D116"; Value=9}
    @{Cell="E116"; Value='[b''document'']'}
    @{Cell="F116"; Value=1}
    @{Cell="D117"; Value=10}
    @{Cell="E117"; Value='[b''back'']'}
    @{Cell="F117"; Value=1}
    @{Cell="D118"; Value=11}
    @{Cell="E118"; Value='[b''in'']'}
    @{Cell="F118"; Value=1}
    @{Cell="D119"; Value=12}
    @{Cell="E119"; Value='[b''and'']'}
    @{Cell="F119"; Value=1}
    @{Cell="D120"; Value=13}
    @{Cell="E120"; Value='[b''has'']'}
    @{Cell="F120"; Value=1}
    @{Cell="D121"; Value=14}
    @{Cell="E121"; Value='[b''played'']'}
    @{Cell="F121"; Value=1}
    @{Cell="D122"; Value=15}
    @{Cell="E122"; Value='[b''a'']'}
    @{Cell="F122"; Value=1}
    @{Cell="D123"; Value=16}
    @{Cell="E123"; Value='[b''strong'']'}
    @{Cell="F123"; Value=1}
    @{Cell="D124"; Value=17}
    @{Cell="E124"; Value='[b''role'']'}
    @{Cell="F124"; Value=1}
    @{Cell="D125"; Value=18}
    @{Cell="E125"; Value='[b''in'']'}
    @{Cell="F125"; Value=1}
    @{Cell="D126"; Value=19}
    @{Cell="E126"; Value='[b''supporting'']'}
    @{Cell="F126"; Value=1}
    @{Cell="C128"; Value=6}
    @{Cell="D128"; Value=0}
    @{Cell="D129"; Value=1}
    @{Cell="D130"; Value=2}
    @{Cell="D131"; Value=3}
    @{Cell="D132"; Value=4}
    @{Cell="D133"; Value=5}
    @{Cell="D134"; Value=6}
    @{Cell="D135"; Value=7}
    @{Cell="D136"; Value=8}
    @{Cell="D137"; Value=9}
    @{Cell="D138"; Value=10}
    @{Cell="D139"; Value=11}
    @{Cell="D140"; Value=12}
    @{Cell="D141"; Value=13}
    @{Cell="D142"; Value=14}
    @{Cell="D143"; Value=15}
    @{Cell="D144"; Value=16}
    @{Cell="D145"; Value=17}
    @{Cell="D146"; Value=18}
    @{Cell="D147"; Value=19}
)
foreach ($item in $cellValues) {
    $ws.Range($item.Cell).Value = $item.Value
}

# --- "padding" cell fills (red = unused word slot, green = sentence-index marker) ---
$cellFills = @(
    @{Cell="D56"; Color=255}
    @{Cell="E56"; Color=255}
    @{Cell="F56"; Color=255}
    @{Cell="D57"; Color=255}
    @{Cell="E57"; Color=255}
    @{Cell="F57"; Color=255}
    @{Cell="D65"; Color=255}
    @{Cell="E65"; Color=255}
    @{Cell="F65"; Color=255}
    @{Cell="D66"; Color=255}
    @{Cell="E66"; Color=255}
    @{Cell="F66"; Color=255}
    @{Cell="D67"; Color=255}
    @{Cell="E67"; Color=255}
    @{Cell="F67"; Color=255}
    @{Cell="D68"; Color=255}
    @{Cell="E68"; Color=255}
    @{Cell="F68"; Color=255}
    @{Cell="C86"; Color=5296274}
    @{Cell="C107"; Color=5296274}
    @{Cell="D128"; Color=255}
    @{Cell="E128"; Color=255}
    @{Cell="F128"; Color=255}
)
foreach ($item in $cellFills) {
    $ws.Range($item.Cell).Interior.Color = $item.Color
}

# --- restore the active selection to match the scrolled-to cell after editing ---
$ws.Range("F134").Select()
